$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 56, pushing existing rows 56:185 down to 57:186
$ws.Range("A56").EntireRow.Insert()

# Populate the newly inserted row 56 with the new weekly price-report record
$ws.Range("A56").Value = 4
$ws.Range("B56").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C56").Value = "Los Lagos"
$ws.Range("D56").Value = 44536
$ws.Range("E56").Value = 10
$ws.Range("F56").Value = 100112003
$ws.Range("G56").Value = "Ajo"
$ws.Range("H56").Value = "Chino"
$ws.Range("I56").Value = "Primera"
$ws.Range("J56").Value = 60
$ws.Range("K56").Value = 21000
$ws.Range("L56").Value = 22000
$ws.Range("M56").Value = 21500
$ws.Range("N56").Value = "$/caja 10 kilos"
$ws.Range("O56").Value = "China"
$ws.Range("P56").Value = 2150
$ws.Range("Q56").Value = 10
$ws.Range("R56").Value = "Hortaliza"
